$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before column N (14), shifting old N:P to O:Q.
# Excel copies the formatting (incl. width) of the column to the left (M) onto
# the freshly inserted column, so mirror that explicitly.
$ws.Columns.Item(14).Insert()
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Make "Repayment Schedule" the active sheet/tab (was "Transactions"),
# and move the selection to P6 (post-insert column layout).
$ws.Activate()
$ws.Range("P6").Select()
